$d = $word.ActiveDocument
$Q = [char]34

# --- Step 1: remove the "Meta description" paragraph near the top of the document ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Meta description*") {
        $p.Range.Delete()
        break
    }
}

# --- Step 2: find the paragraph that still holds the old "Create a cartoon-style..." image
#             prompt (the last paragraph in the document) and insert a new bold title
#             paragraph right before it. ---
$oldImagePromptStart = "Create a cartoon-style feature image for Dragons of the North"
$imagePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "$oldImagePromptStart*") {
        $imagePara = $p
        break
    }
}

$newPara = $imagePara.Range.InsertParagraphBefore()

# Re-locate the freshly inserted (now empty) paragraph and give it the exact target
# run/formatting structure (leading empty run + bold run with the title text) via a
# WordprocessingML fragment, mirroring the sibling paragraphs in this document.
$titlePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "$oldImagePromptStart*") {
        $titlePara = $d.Paragraphs.Item($i - 1)
        break
    }
}

$titleXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Dragons of the North Free Slot Game | Pro Review</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$titlePara.Range.InsertXML($titleXml)

# --- Step 3: replace the old italic image-prompt text with the meta-description body copy,
#             keeping the paragraph's existing italic run formatting intact. ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "$oldImagePromptStart*") {
        [void]$p.Range.Find.Execute(
            "Create a cartoon-style feature image for Dragons of the North that showcases a happy Maya warrior with glasses. The Maya warrior should be riding a dragon and holding a dragon egg in one hand. The background should feature medieval castles, fire, and other dragons flying in the sky. The image should be vibrant and exciting to attract players to the game. The text " + $Q + "Dragons of the North" + $Q + " should be prominently displayed in a fun font that matches the theme of the game. The image should be suitable for use in online advertisements, social media posts, and other marketing materials.",
            $true, $false, $false, $false, $false, $true, 1, $false,
            "Read our unbiased review of Dragons of the North, a 5-reel, 50-payline slot game with various bonuses and jackpots. Play it for free and discover its pros and cons.",
            2)
        break
    }
}

Write-Host "Moved the meta description into the closing Play/Read paragraphs."
